$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Update the phone number value in C3 (last digit 1 -> 0)
$ws.Range("C3").Value = 6923378500

# Remove rows 4 through 11 (extra generated records no longer needed)
$ws.Range("A4:C11").EntireRow.Delete()

# Update the selected cell to reflect the new extent of data
$ws.Range("C3").Select()
